$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Select()

# --- New shared-string entries / rows -------------------------------------

# "Mirror" entry (row 43)
$ws.Range("A43").Value = "Mirror"
$ws.Range("G43").Value = "Focal Length: 1.6m"

# "Arduino" entry (rows 45-46)
$ws.Range("A45").Value = "Arduino"
$ws.Range("B45").Value = "ZYDUINO UNO"
$ws.Range("C45").Value = "5V for operating voltage"
$ws.Range("C46").Value = "7-12V for input voltage"

# --- Update the window/view state to match the final workbook -------------

$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$win.Zoom = 113
$ws.Range("C46").Select()
